# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig
# Updates the StructureDefinition metadata (URL, Version, Date, Publisher)
# on the "Metadata" sheet, and clears the now-redundant constraint text
# duplicated on the base "Extension" row of the "Elements" sheet (the
# same ele-1/ext-1 constraint text remains on the Extension.extension row).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/job-class-code"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
# A leading single-quote forces the cell to stay text-typed with an empty
# string (rather than Value="" which clears the cell to a blank/Number
# cell) -- matches the source file, where this cell still stores an
# (empty) shared string rather than being truly blank.
$elements.Range("AI2").Value = "'"

# Extension.url's "Fixed Value" column echoes the same URL string as
# Metadata!B2 in the source document; keep it in sync.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/job-class-code"
